$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1875
$ws.Range("C2").Value = 0.5625
$ws.Range("P2").Value = 0.125
$ws.Range("S2").Value = 0.125
$ws.Range("P3").Value = 0.4444444444444444
$ws.Range("S3").Value = 0.5555555555555556
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.2
$ws.Range("F6").Value = 0.1
$ws.Range("J6").Value = 0.05
$ws.Range("R6").Value = 0.05
$ws.Range("S6").Value = 0.6
$ws.Range("B7").Value = 0.06666666666666667
$ws.Range("F7").Value = 0.1333333333333333
$ws.Range("J7").Value = 0.06666666666666667
$ws.Range("Q7").Value = 0.1333333333333333
$ws.Range("R7").Value = 0.06666666666666667
$ws.Range("S7").Value = 0.5333333333333333
$ws.Range("B8").Value = 0.09090909090909091
$ws.Range("D8").Value = 0.09090909090909091
$ws.Range("J8").Value = 0.1515151515151515
$ws.Range("O8").Value = 0.0303030303030303
$ws.Range("Q8").Value = 0.09090909090909091
$ws.Range("R8").Value = 0.1212121212121212
$ws.Range("S8").Value = 0.3333333333333333
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("E9").Value = 0.04545454545454546
$ws.Range("F9").Value = 0.04545454545454546
$ws.Range("J9").Value = 0.1363636363636364
$ws.Range("Q9").Value = 0.04545454545454546
$ws.Range("R9").Value = 0.1363636363636364
$ws.Range("S9").Value = 0.5
$ws.Range("B10").Value = 0.04054054054054054
$ws.Range("D10").Value = 0.01351351351351351
$ws.Range("F10").Value = 0.1216216216216216
$ws.Range("J10").Value = 0.1081081081081081
$ws.Range("O10").Value = 0.04054054054054054
$ws.Range("Q10").Value = 0.1486486486486487
$ws.Range("R10").Value = 0.06756756756756757
$ws.Range("S10").Value = 0.4594594594594595
$ws.Range("K11").Value = 0.25
$ws.Range("S11").Value = 0.05
$ws.Range("G12").Value = 0.9090909090909091
$ws.Range("J12").Value = 0.09090909090909091
$ws.Range("H15").Value = 0.06666666666666667
$ws.Range("I15").Value = 0.1333333333333333
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.4
$ws.Range("H16").Value = 0.1111111111111111
$ws.Range("I16").Value = 0.2222222222222222
$ws.Range("J16").Value = 0.2222222222222222
$ws.Range("O16").Value = 0.2222222222222222
$ws.Range("S16").Value = 0.2222222222222222
$ws.Range("H17").Value = 0.2941176470588235
$ws.Range("I17").Value = 0.1764705882352941
$ws.Range("J17").Value = 0.1764705882352941
$ws.Range("K17").Value = 0.05882352941176471
$ws.Range("O17").Value = 0.1176470588235294
$ws.Range("S17").Value = 0.1764705882352941
$ws.Range("H18").Value = 0.07692307692307693
$ws.Range("I18").Value = 0.2307692307692308
$ws.Range("J18").Value = 0.3846153846153846
$ws.Range("O18").Value = 0.07692307692307693
$ws.Range("S18").Value = 0.2307692307692308
$ws.Range("F19").Value = 0.01754385964912281
$ws.Range("H19").Value = 0.2280701754385965
$ws.Range("I19").Value = 0.1052631578947368
$ws.Range("J19").Value = 0.3596491228070176
$ws.Range("K19").Value = 0.1052631578947368
$ws.Range("M19").Value = 0.01754385964912281
$ws.Range("O19").Value = 0.03508771929824561
$ws.Range("S19").Value = 0.131578947368421
